$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$sub5 = [char]8325
$sub8 = [char]8328

# Column D (Price) updates - force text format so numeric-looking strings
# (e.g. "214.25", "1.001") are preserved exactly as text, matching the
# original inline-string cell contents.
$dCells = @{
    "D2"  = "26.089.38"
    "D3"  = "1.641.40"
    "D5"  = "214.25"
    "D6"  = "0.5228"
    "D8"  = "0.2604"
    "D9"  = "0.06307"
    "D10" = "20.66"
    "D11" = "0.07672"
    "D12" = "1.649.66"
    "D13" = "4.411"
    "D14" = "1.859.90"
    "D15" = "0.5535"
    "D16" = "0.0${sub5}8191"
    "D17" = "65.02"
    "D18" = "26.086.87"
    "D19" = "1.001"
    "D20" = "4.703"
    "D21" = "188.77"
    "D22" = "10.20"
    "D23" = "6.171"
    "D24" = "1.001"
    "D25" = "145.27"
    "D27" = "7.407"
    "D28" = "15.81"
    "D29" = "1.385"
    "D30" = "0.05941"
    "D31" = "1.257"
    "D32" = "3.438"
    "D33" = "3.400"
    "D34" = "1.650"
    "D35" = "0.9830"
    "D36" = "2.391"
    "D37" = "2.755"
    "D38" = "0.5632"
    "D39" = "0.01617"
    "D40" = "0.8518"
    "D42" = "5.692"
    "D43" = "1.022.82"
    "D44" = "100.15"
    "D45" = "1.788.61"
    "D46" = "0.0${sub8}108"
    "D47" = "55.80"
    "D48" = "1.001"
    "D49" = "8.079"
    "D50" = "0.05148"
    "D51" = "0.4217"
}

foreach ($addr in $dCells.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $dCells[$addr]
}

# Column E (Volume(1h)) updates
$eCells = @{
    "E2"  = "  -0.15%  "
    "E3"  = "  -1.56%  "
    "E4"  = "  -0.15%  "
    "E5"  = "  +2.37%  "
    "E6"  = "  +0.02%  "
    "E7"  = "  -0.16%  "
    "E8"  = "  -0.30%  "
    "E9"  = "  -0.34%  "
    "E10" = "  -2.07%  "
    "E11" = "  +1.94%  "
    "E12" = "  -1.12%  "
    "E13" = "  -0.37%  "
    "E14" = "  -1.86%  "
    "E15" = "  +0.77%  "
    "E16" = "  +2.95%  "
    "E17" = "  -2.13%  "
    "E18" = "  -0.19%  "
    "E19" = "  -0.15%  "
    "E20" = "  -0.07%  "
    "E21" = "  +1.12%  "
    "E22" = "  -0.45%  "
    "E23" = "  -0.14%  "
    "E24" = "  -0.19%  "
    "E25" = "  -2.94%  "
    "E26" = "  -3.15%  "
    "E27" = "  -0.93%  "
    "E28" = "  -0.23%  "
    "E29" = "  +2.96%  "
    "E30" = "  -6.32%  "
    "E31" = "  -1.20%  "
    "E32" = "  -1.29%  "
    "E33" = "  -0.29%  "
    "E34" = "  +0.77%  "
    "E35" = "  -1.99%  "
    "E36" = "  -0.67%  "
    "E37" = "  +0.48%  "
    "E38" = "  -6.17%  "
    "E39" = "  +0.25%  "
    "E40" = "  -1.73%  "
    "E41" = "  -0.23%  "
    "E42" = "  -6.82%  "
    "E43" = "  -7.80%  "
    "E44" = "  +0.25%  "
    "E45" = "  -1.61%  "
    "E46" = "  -1.51%  "
    "E48" = "  +0.21%  "
    "E49" = "  +0.66%  "
    "E50" = "  -1.54%  "
    "E51" = "  -0.63%  "
}

foreach ($addr in $eCells.Keys) {
    $ws.Range($addr).Value = $eCells[$addr]
}
